$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add New York state hospitalization data for 14 April 2020 (row 31)
$ws.Range("A31").Value = 43935
$ws.Range("B31").Value = -362
$ws.Range("C31").Value = -20
$ws.Range("D31").Value = -7
$ws.Range("F31").Value = 752
$ws.Range("G31").Value = 2253

# Copy the date number format from the row above (A30) so A31 renders as YYYY-MM-DD
$ws.Range("A30").Copy()
$ws.Range("A31").PasteSpecial(-4122)

# Match the selection state from the diff (active cell moved to F32)
$ws.Range("F32").Select()
